$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "42.296.20"
$ws.Range("E2").Value = "  +0.72%  "

# Row 3
Set-TextValue "D3" "2.290.32"
$ws.Range("E3").Value = "  -0.43%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
Set-TextValue "D5" "321.46"
$ws.Range("E5").Value = "  +1.43%  "

# Row 6
Set-TextValue "D6" "102.71"
$ws.Range("E6").Value = "  -1.36%  "

# Row 7
$ws.Range("E7").Value = "  -0.81%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -0.36%  "

# Row 10
Set-TextValue "D10" "39.63"
$ws.Range("E10").Value = "  -0.71%  "

# Row 11
Set-TextValue "D11" "0.0904"
$ws.Range("E11").Value = "  -0.77%  "

# Row 12
Set-TextValue "D12" "8.31"
$ws.Range("E12").Value = "  -1.85%  "

# Row 13
$ws.Range("E13").Value = "  -0.98%  "

# Row 14
Set-TextValue "D14" "0.964"
$ws.Range("E14").Value = "  -1.39%  "

# Row 15
Set-TextValue "D15" "15.16"
$ws.Range("E15").Value = "  -1.64%  "

# Row 16
Set-TextValue "D16" "2.637.54"
$ws.Range("E16").Value = "  -0.42%  "

# Row 17
Set-TextValue "D17" "2.290.36"
$ws.Range("E17").Value = "  -0.15%  "

# Row 18
Set-TextValue "D18" "42.281.13"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19
$ws.Range("E19").Value = "  -4.26%  "

# Row 20
$ws.Range("E20").Value = "  -0.43%  "

# Row 21
Set-TextValue "D21" "12.68"
$ws.Range("E21").Value = "  +27.89%  "

# Row 22
Set-TextValue "D22" "73.02"
$ws.Range("E22").Value = "  -0.17%  "

# Row 23
$ws.Range("E23").Value = "  -0.46%  "

# Row 24
Set-TextValue "D24" "268.49"
$ws.Range("E24").Value = "  +3.48%  "

# Row 25
Set-TextValue "D25" "2.23"
$ws.Range("E25").Value = "  -3.50%  "

# Row 26
$ws.Range("E26").Value = "  -0.25%  "

# Row 27
Set-TextValue "D27" "10.86"
$ws.Range("E27").Value = "  -0.66%  "

# Row 28
Set-TextValue "D28" "2.35"
$ws.Range("E28").Value = "  +5.04%  "

# Row 29
Set-TextValue "D29" "22.55"
$ws.Range("E29").Value = "  -2.80%  "

# Row 30
Set-TextValue "D30" "38.03"
$ws.Range("E30").Value = "  +5.87%  "

# Row 31
Set-TextValue "D31" "164.48"
$ws.Range("E31").Value = "  +0.19%  "

# Row 32
$ws.Range("E32").Value = "  +2.93%  "

# Row 33
Set-TextValue "D33" "0.0872"
$ws.Range("E33").Value = "  -1.73%  "

# Row 34
$ws.Range("E34").Value = "  +1.00%  "

# Row 35
$ws.Range("E35").Value = "  -5.46%  "

# Row 36
Set-TextValue "D36" "2.50"
$ws.Range("E36").Value = "  -14.18%  "

# Row 37
Set-TextValue "D37" "4.59"
$ws.Range("E37").Value = "  -1.24%  "

# Row 38
$ws.Range("E38").Value = "  +0.96%  "

# Row 39
Set-TextValue "D39" "3.69"
$ws.Range("E39").Value = "  +0.52%  "

# Row 40
$ws.Range("E40").Value = "  -6.07%  "

# Row 41
$ws.Range("E41").Value = "  +3.48%  "

# Row 42: was FirstDigitalUSD -> now MultiversX (rows 42/43 swap identities)
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D42" "68.72"
$ws.Range("E42").Value = "  -3.30%  "

# Row 43: was MultiversX -> now FirstDigitalUSD
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "1.00"
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("E44").Value = "  -0.84%  "

# Row 45
Set-TextValue "D45" "90.88"
$ws.Range("E45").Value = "  -10.28%  "

# Row 46
$ws.Range("E46").Value = "  +1.24%  "

# Row 47
Set-TextValue "D47" "113.88"
$ws.Range("E47").Value = "  -0.42%  "

# Row 48
Set-TextValue "D48" "80.02"
$ws.Range("E48").Value = "  +1.44%  "

# Row 49
Set-TextValue "D49" "8.97"
$ws.Range("E49").Value = "  -0.94%  "

# Row 50
Set-TextValue "D50" "1.606.01"
$ws.Range("E50").Value = "  +3.90%  "

# Row 51
$ws.Range("E51").Value = "  -1.87%  "

